$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4651.067
$ws.Range("I40").Value = 7240
$ws.Range("K40").Value = 7240
$ws.Range("M40").Value = -7065

$ws.Range("H51").Value = 8548.25
$ws.Range("J51").Value = 8508.091
$ws.Range("L51").Value = 8508.091
$ws.Range("N51").Value = -9476.091

$ws.Range("I64").Value = 2883.125
$ws.Range("J64").Value = 3999
$ws.Range("K64").Value = 2883.125
$ws.Range("L64").Value = 3999
$ws.Range("M64").Value = -2635.125
$ws.Range("N64").Value = -4495

$ws.Range("I67").Value = 2883.125
$ws.Range("J67").Value = 3999
$ws.Range("K67").Value = 2883.125
$ws.Range("L67").Value = 3999
$ws.Range("M67").Value = -2025.125
$ws.Range("N67").Value = -5715

$ws.Range("H76").Value = 3415.3333
$ws.Range("I76").Value = 3415.3333
$ws.Range("K76").Value = 3415.3333
$ws.Range("M76").Value = -3100.3333

$ws.Range("H79").Value = 3415.3333
$ws.Range("I79").Value = 3415.3333
$ws.Range("K79").Value = 3415.3333
$ws.Range("M79").Value = -2323.3333

$ws.Range("H130").Value = 55800
$ws.Range("J130").Value = 55800
$ws.Range("L130").Value = 55800
$ws.Range("N130").Value = -65840

$ws.Range("H137").Value = 63264.668
$ws.Range("I137").Value = 79555.14
$ws.Range("J137").Value = 6248
$ws.Range("K137").Value = 238665.42
$ws.Range("L137").Value = 18744
$ws.Range("M137").Value = -236115.42
$ws.Range("N137").Value = -23844

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 709
$ws.Range("I30").Value = 709
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 709
$ws.Range("L30").Value = 0
$ws.Range("M30").ClearContents()
$ws.Range("N30").Value = -559

$ws.Range("H74").Value = 52502.49
$ws.Range("I74").Value = 64598.324
$ws.Range("J74").Value = 5631.125
$ws.Range("K74").Value = 64598.324
$ws.Range("L74").Value = 5631.125
$ws.Range("M74").Value = -63724.324
$ws.Range("N74").Value = -7379.125

$ws.Range("H77").Value = 52502.49
$ws.Range("I77").Value = 64598.324
$ws.Range("J77").Value = 5631.125
$ws.Range("K77").Value = 322991.62
$ws.Range("L77").Value = 28155.625
$ws.Range("M77").Value = -318623.62
$ws.Range("N77").Value = -36891.625

$ws.Range("H102").Value = 2296.6667
$ws.Range("I102").Value = 1995.7273
$ws.Range("K102").Value = 1995.7273
$ws.Range("M102").Value = -373.7273

$ws.Range("H122").Value = 4538.86
$ws.Range("I122").Value = 4651.4863
$ws.Range("J122").Value = 4218.3076
$ws.Range("K122").Value = 13954.4589
$ws.Range("L122").Value = 12654.9228
$ws.Range("M122").Value = -11504.4589
$ws.Range("N122").Value = -17554.9228

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3171.1345
$ws.Range("I134").Value = 3217
$ws.Range("K134").Value = 9651
$ws.Range("M134").Value = -7116

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 29708566
$ws.Range("J31").Value = 76926520
$ws.Range("L31").Value = 76926520
$ws.Range("N31").Value = -76927110

$ws.Range("H34").Value = 29708566
$ws.Range("J34").Value = 76926520
$ws.Range("L34").Value = 76926520
$ws.Range("N34").Value = -76926924

$ws.Range("H58").Value = 2412
$ws.Range("J58").Value = 4830
$ws.Range("L58").Value = 4830
$ws.Range("N58").Value = -5236

$ws.Range("H59").Value = 57157.715
$ws.Range("I59").Value = 50052
$ws.Range("K59").Value = 50052
$ws.Range("M59").Value = -48907

$ws.Range("H122").Value = 2177.2
$ws.Range("I122").Value = 1967.7142
$ws.Range("K122").Value = 5903.142599999999
$ws.Range("M122").Value = -3453.142599999999

$ws.Range("H134").Value = 4907.222
$ws.Range("I134").Value = 3784.75
$ws.Range("K134").Value = 11354.25
$ws.Range("M134").Value = -8819.25

$ws.Range("H136").Value = 2412
$ws.Range("J136").Value = 4830
$ws.Range("L136").Value = 14490
$ws.Range("N136").Value = -19590

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 47622904
$ws.Range("I64").Value = 3999.5
$ws.Range("J64").Value = 66670468
$ws.Range("K64").Value = 11998.5
$ws.Range("L64").Value = 200011404
$ws.Range("M64").Value = -11728.5
$ws.Range("N64").Value = -200011944

$ws.Range("H67").Value = 47622904
$ws.Range("I67").Value = 3999.5
$ws.Range("J67").Value = 66670468
$ws.Range("K67").Value = 11998.5
$ws.Range("L67").Value = 200011404
$ws.Range("M67").Value = -11062.5
$ws.Range("N67").Value = -200013276

$ws.Range("H68").Value = 5557681
$ws.Range("J68").Value = 2589.2727
$ws.Range("L68").Value = 7767.8181
$ws.Range("N68").Value = -9389.8181

$ws.Range("H71").Value = 5557681
$ws.Range("J71").Value = 2589.2727
$ws.Range("L71").Value = 23303.4543
$ws.Range("N71").Value = -31415.4543

$ws.Range("H132").Value = 7345.05
$ws.Range("I132").Value = 10065.846
$ws.Range("J132").Value = 2292.1428
$ws.Range("K132").Value = 90592.614
$ws.Range("L132").Value = 20629.2852
$ws.Range("M132").Value = -88062.614
$ws.Range("N132").Value = -25689.2852

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 9575
$ws.Range("I55").Value = 9575
$ws.Range("K55").Value = 9575
$ws.Range("M55").Value = -9248

$ws.Range("H102").Value = 24728.4
$ws.Range("I102").Value = 2164.6943
$ws.Range("K102").Value = 2164.6943
$ws.Range("M102").Value = -542.6943000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4308.091
$ws.Range("I40").Value = 4308.091
$ws.Range("K40").Value = 4308.091
$ws.Range("M40").Value = -4172.091

$ws.Range("H55").Value = 6464.1924
$ws.Range("I55").Value = 751.2222
$ws.Range("J55").Value = 19318.375
$ws.Range("K55").Value = 751.2222
$ws.Range("L55").Value = 19318.375
$ws.Range("M55").Value = -578.2222
$ws.Range("N55").Value = -19664.375

$ws.Range("H122").Value = 2877.6
$ws.Range("J122").Value = 5000
$ws.Range("L122").Value = 15000
$ws.Range("N122").Value = -19900

$ws.Range("H136").Value = 7985.25
$ws.Range("I136").Value = 8554.714
$ws.Range("K136").Value = 25664.142
$ws.Range("M136").Value = -23114.142

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1226.6666
$ws.Range("I100").Value = 1540
$ws.Range("J100").Value = 600
$ws.Range("K100").Value = 3080
$ws.Range("L100").Value = 1200
$ws.Range("M100").Value = -2539
$ws.Range("N100").Value = -2282

$ws.Range("H126").Value = 2126.923
$ws.Range("I126").Value = 2160.8696
$ws.Range("J126").Value = 1866.6666
$ws.Range("K126").Value = 6482.6088
$ws.Range("L126").Value = 5599.9998
$ws.Range("M126").Value = -4012.6088
$ws.Range("N126").Value = -10539.9998

$ws.Range("H136").Value = 419590
$ws.Range("I136").Value = 419590
$ws.Range("K136").Value = 1258770
$ws.Range("M136").Value = -1256220

Write-Output "Applied all Lich Profits updates"